$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("A2").ClearContents()
$ws.Range("B2").Value = " <color=#00CC00>(We’ve returned once again to the place where the Lord’s body was discovered last night.)</color>"
$ws.Range("D2").Value = "DialogueVocal"
$ws.Range("E2").Value = "Door-Investigate"
$ws.Range("F2").Value = "Suspicious"
$ws.Range("G2:I2").Clear()
$ws.Range("J2").Value = "appearAt"
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = "Dee-Thinking"
$ws.Rows.Item(2).RowHeight = 51

# --- Row 3 ---
$ws.Range("A3").ClearContents()
$ws.Range("B3").Value = " <color=#00CC00>(The scene remains exactly as it was—unchanged from yesterday.)</color>"
$ws.Range("D3").Value = "DialogueVocal"
$ws.Range("E3").Value = "Door-Investigate"
$ws.Range("G3:I3").Clear()
$ws.Rows.Item(3).RowHeight = 34

# --- Row 4 ---
$ws.Range("A4").ClearContents()
$ws.Range("B4").Value = " <color=#00CC00>(But under the bright daylight, we’re bound to uncover new clues.)</color>"
$ws.Range("D4").Value = "DialogueVocal"
$ws.Range("E4").Value = "Door-Investigate"
$ws.Range("G4:I4").Clear()
$ws.Rows.Item(4).RowHeight = 34

# --- Row 5 ---
$ws.Range("A5").Value = "Dee"
$ws.Range("B5").Value = "Let’s begin the investigation!"
$ws.Range("C5").Value = "Dee-Regular"
$ws.Range("D5").Value = "DialogueVocal"
$ws.Range("E5").Value = "Door-Investigate"
$ws.Range("G5:I5").Clear()

# --- Row 6 ---
$ws.Range("A6").Value = "Investigate2"
$ws.Range("B6").Value = "Water"
$ws.Range("C6").Value = "Water"
$ws.Range("D6").Value = "DialogueVocal"
$ws.Range("E6").Value = "Door-Investigate"
$ws.Range("G6:I6").Clear()

# --- Row 7 ---
$ws.Range("B7").Value = "Hand"
$ws.Range("C7").Value = "Hand"
$ws.Range("D7").Value = "DialogueVocal"
$ws.Range("E7").Value = "Door-Investigate"
$ws.Range("G7:I7").Clear()

# --- Row 8 ---
$ws.Range("B8").Value = "Blood"
$ws.Range("C8").Value = "Blood"
$ws.Range("D8").Value = "DialogueVocal"
$ws.Range("E8").Value = "Door-Investigate"
$ws.Range("G8:I8").Clear()

# --- Row 9 (new data row, was blank formatting-only row before) ---
$ws.Range("B9").Value = "End Investigation"
$ws.Range("C9").Value = "StoryScript14"
$ws.Range("D9").Value = "DialogueVocal"
$ws.Range("E9").Value = "Door-Investigate"
$ws.Range("G9:I9").Clear()
$ws.Range("J9").Value = "disappear"
$ws.Rows.Item(9).RowHeight = 17

# --- Remove now-unused rows 10:16 ---
$ws.Range("A10:P16").EntireRow.Delete()

# --- Selection ---
$ws.Range("B7").Select()

Write-Host "Edit complete"
